# Update EC (Estado de Cuenta) database: remove JOSE DAVID PEÑA CABARCAS entry,
# regroup/update the remaining worker rows, and refresh the summary totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1) Remove the row for JOSE DAVID PEÑA CABARCAS (period 2507). Everything below
#    shifts up by one, so SANTIAGO's row (which already carries the bottom-border
#    style for the last row of the table) becomes the new last data row.
$ws.Rows("22").Delete()

# 2) Rewrite the worker/period rows with the updated data set, grouped by worker.
$ws.Range("C16").Value = "9154405"
$ws.Range("D16").Value = "WILFRAN ANTONIO SALGADO MORALES"
$ws.Range("E16").Value = "2107"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 908526

$ws.Range("C17").Value = "9158750"
$ws.Range("D17").Value = "NESTOR CASSIANI BELLO"
$ws.Range("E17").Value = "2107"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 908526

$ws.Range("C18").Value = "9154405"
$ws.Range("D18").Value = "WILFRAN ANTONIO SALGADO MORALES"
$ws.Range("E18").Value = "2108"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 908526

$ws.Range("C19").Value = "9158750"
$ws.Range("D19").Value = "NESTOR CASSIANI BELLO"
$ws.Range("E19").Value = "2108"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 908526

$ws.Range("C20").Value = "45580898"
$ws.Range("D20").Value = "SARA CRISTINA ROCHA MARTELO"
$ws.Range("E20").Value = "2304"
$ws.Range("F20").Value = 9280
$ws.Range("G20").Value = 1423500

$ws.Range("C21").Value = "45580898"
$ws.Range("D21").Value = "SARA CRISTINA ROCHA MARTELO"
$ws.Range("E21").Value = "2305"
$ws.Range("F21").Value = 23200
$ws.Range("G21").Value = 1423500

$ws.Range("C22").Value = "1044906536"
$ws.Range("D22").Value = "SANTIAGO ELIAS VARGAS CASTRO"
$ws.Range("E22").Value = "2401"
$ws.Range("F22").Value = 10400
$ws.Range("G22").Value = 1423500

# 3) Refresh the summary totals at the top of the sheet.
$ws.Range("E11").Value = 188244
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 5
